# Reorder/update the menu items in rows 5-55 of Sheet1 (A/C/D columns)
# to reflect the new grouping of the restaurant menu (pizzas, sandwiches,
# snacks, bakery & beverages moved/re-sorted).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(5, 1).Value = '7 Inch Pizza Margarita'
$ws.Cells.Item(5, 3).Value = 90
$ws.Cells.Item(5, 4).Value = '7 Inch Margrita Pizza.png'
$ws.Cells.Item(6, 1).Value = '7 Inch Pizza Onion and Capsicum'
$ws.Cells.Item(6, 3).Value = 90
$ws.Cells.Item(6, 4).Value = '7 Inch Onion and Capsicum Pizza.png'
$ws.Cells.Item(7, 1).Value = '7 Inch Pizza Paneer'
$ws.Cells.Item(7, 3).Value = 90
$ws.Cells.Item(7, 4).Value = '7 Inch Paneer Pizza.png'
$ws.Cells.Item(8, 1).Value = '7 Inch Pizza Corn'
$ws.Cells.Item(8, 3).Value = 90
$ws.Cells.Item(8, 4).Value = '7 Inch Corn Pizza.png'
$ws.Cells.Item(9, 1).Value = '7 Inch Combo Pack (Pizza Margarita, Onion and Capsicum, Corn, Pizza Paneer)'
$ws.Cells.Item(9, 3).Value = 350
$ws.Cells.Item(9, 4).Value = '7 Inch Combo Pack.png'
$ws.Cells.Item(10, 1).Value = 'Pizza Margarita 10'' Inch'
$ws.Cells.Item(10, 3).Value = 150
$ws.Cells.Item(10, 4).Value = '10 Inch Margrita Pizza.png'
$ws.Cells.Item(11, 1).Value = 'Pizza Corn 10'' Inch'
$ws.Cells.Item(11, 3).Value = 150
$ws.Cells.Item(11, 4).Value = '10 Inch Corn Pizza.png'
$ws.Cells.Item(12, 1).Value = 'Pizza Onion and Capsicum 10'' Inch'
$ws.Cells.Item(12, 3).Value = 175
$ws.Cells.Item(12, 4).Value = '10 Inch Onion and Capsicum Pizza.png'
$ws.Cells.Item(13, 1).Value = 'Pizza Paneer, Veggie ( Onion and Capsicum and corn) 10'' Inch'
$ws.Cells.Item(13, 3).Value = 200
$ws.Cells.Item(13, 4).Value = '10 Inch Pizza Veggi Panner( Onion and Capcium and corn).png'
$ws.Cells.Item(14, 1).Value = 'Tortila Wrap Paneer'
$ws.Cells.Item(14, 3).Value = 100
$ws.Cells.Item(14, 4).Value = 'Tortila Wrap Paneer.jpeg'
$ws.Cells.Item(15, 1).Value = 'Burrito wrap Paneer'
$ws.Cells.Item(15, 3).Value = 100
$ws.Cells.Item(15, 4).Value = 'Burrito wrap Paneer.jpeg'
$ws.Cells.Item(16, 1).Value = 'Quesadilla Paneer'
$ws.Cells.Item(16, 3).Value = 100
$ws.Cells.Item(16, 4).Value = 'Quesadilla Paneer.jpeg'
$ws.Cells.Item(17, 1).Value = 'Schezwan Grilled Sandwich – Indo-Chinese fusion with Schezwan sauce, Veggies Patty, and cheese.'
$ws.Cells.Item(17, 3).Value = 100
$ws.Cells.Item(17, 4).Value = 'Schezwan Grilled Sandwich.png'
$ws.Cells.Item(18, 1).Value = 'Club Sandwich (Indian Style) –Paneer Patty, and cheese, Veggies, green chutney'
$ws.Cells.Item(18, 3).Value = 100
$ws.Cells.Item(18, 4).Value = 'Multi-layered with veggies Paneer.png'
$ws.Cells.Item(19, 1).Value = 'Veg Biryani Soya with Garlic Mayo Dip *1'
$ws.Cells.Item(19, 3).Value = 120
$ws.Cells.Item(19, 4).Value = 'Veg Biryani Soya with Garlic Mayo Dip 1.png'
$ws.Cells.Item(20, 1).Value = 'Veg Biryani Paneer with Dip Garlic Mayo Dip *1'
$ws.Cells.Item(20, 3).Value = 175
$ws.Cells.Item(20, 4).Value = 'Veg Paneer Biryani.png'
$ws.Cells.Item(21, 1).Value = 'Salted French Frise'
$ws.Cells.Item(21, 3).Value = 80
$ws.Cells.Item(21, 4).Value = 'Salted French Fires image.jpg'
$ws.Cells.Item(22, 1).Value = 'Peri-Peri French Frise'
$ws.Cells.Item(22, 3).Value = 100
$ws.Cells.Item(22, 4).Value = 'Peri- Peri French Fires.jpg'
$ws.Cells.Item(23, 1).Value = 'vegetable-mayonnaise-sandwich'
$ws.Cells.Item(23, 3).Value = 20
$ws.Cells.Item(23, 4).Value = 'vegetable-mayonnaise-sandwich.jpg'
$ws.Cells.Item(24, 1).Value = 'Aloo Patty'
$ws.Cells.Item(24, 3).Value = 20
$ws.Cells.Item(24, 4).Value = 'Aloo Patty.jpg'
$ws.Cells.Item(25, 1).Value = 'Paneer Patty'
$ws.Cells.Item(25, 3).Value = 25
$ws.Cells.Item(25, 4).Value = 'Paneer Patty.jpg'
$ws.Cells.Item(26, 1).Value = 'Butter Patty'
$ws.Cells.Item(26, 3).Value = 30
$ws.Cells.Item(26, 4).Value = 'Butter Patty.jpg'
$ws.Cells.Item(27, 1).Value = 'Pastry Pineapple'
$ws.Cells.Item(27, 3).Value = 25
$ws.Cells.Item(27, 4).Value = 'Pastry Pineapple.jpg'
$ws.Cells.Item(28, 1).Value = 'Pastry Chocolate'
$ws.Cells.Item(28, 3).Value = 45
$ws.Cells.Item(28, 4).Value = 'Pastry Chocolate.jpg'
$ws.Cells.Item(29, 1).Value = 'Cake 1 kg Chocolate'
$ws.Cells.Item(29, 3).Value = 550
$ws.Cells.Item(29, 4).Value = 'Cake Chocolate.jpeg'
$ws.Cells.Item(30, 1).Value = 'Cake 1 kg Pineapple'
$ws.Cells.Item(30, 3).Value = 480
$ws.Cells.Item(30, 4).Value = 'Cake Pineapple.jpg'
$ws.Cells.Item(31, 1).Value = 'Cake 1kg Butterscotch'
$ws.Cells.Item(31, 3).Value = 500
$ws.Cells.Item(31, 4).Value = 'Cake Butterscotch.jpg'
$ws.Cells.Item(32, 1).Value = 'Cake 600gm Pineapple'
$ws.Cells.Item(32, 3).Value = 350
$ws.Cells.Item(32, 4).Value = 'Cake Pineapple.jpg'
$ws.Cells.Item(33, 1).Value = 'Cake 600gm Chocolate'
$ws.Cells.Item(33, 3).Value = 450
$ws.Cells.Item(33, 4).Value = 'Cake Chocolate.jpeg'
$ws.Cells.Item(34, 1).Value = 'Cake 600gm Butterscotch'
$ws.Cells.Item(34, 3).Value = 400
$ws.Cells.Item(34, 4).Value = 'Cake Butterscotch.jpg'
$ws.Cells.Item(35, 1).Value = 'Chips 05'
$ws.Cells.Item(35, 3).Value = 5
$ws.Cells.Item(35, 4).Value = 'Chips 5.jpg'
$ws.Cells.Item(36, 1).Value = 'Chips 10'
$ws.Cells.Item(36, 3).Value = 10
$ws.Cells.Item(36, 4).Value = 'Chips 10.jpg'
$ws.Cells.Item(37, 1).Value = 'Chips 20'
$ws.Cells.Item(37, 3).Value = 20
$ws.Cells.Item(37, 4).Value = 'Chips 20.jpg'
$ws.Cells.Item(38, 1).Value = 'Ti Tac 5'
$ws.Cells.Item(38, 3).Value = 5
$ws.Cells.Item(38, 4).Value = 'Tictac 5.jpg'
$ws.Cells.Item(39, 1).Value = 'Bourbon and Dark Fantasy Biscut 10'
$ws.Cells.Item(39, 3).Value = 10
$ws.Cells.Item(39, 4).Value = 'Bourbon and Dark Fantasy Biscut 10.jpg'
$ws.Cells.Item(40, 1).Value = 'Cake 15'
$ws.Cells.Item(40, 3).Value = 15
$ws.Cells.Item(40, 4).Value = 'britannia cake 15.jpg'
$ws.Cells.Item(41, 1).Value = 'Hide and Seek Black bourbon 10'
$ws.Cells.Item(41, 3).Value = 10
$ws.Cells.Item(41, 4).Value = 'Hide and Seek Black bourne 10.jpg'
$ws.Cells.Item(42, 1).Value = 'Hide and Seek 10'
$ws.Cells.Item(42, 3).Value = 10
$ws.Cells.Item(42, 4).Value = 'Hide and seek 10.jpg'
$ws.Cells.Item(43, 1).Value = 'Hide and Seek 30'
$ws.Cells.Item(43, 3).Value = 30
$ws.Cells.Item(43, 4).Value = 'Hide and Seek Biscut 30.jpg'
$ws.Cells.Item(44, 1).Value = '5 Star 5 rs'
$ws.Cells.Item(44, 3).Value = 5
$ws.Cells.Item(44, 4).Value = '5 Star 5 rs.jpg'
$ws.Cells.Item(45, 1).Value = 'dairy milk chocolate 20 rs'
$ws.Cells.Item(45, 3).Value = 20
$ws.Cells.Item(45, 4).Value = 'dairy milk chocolate 20 rs.jpg'
$ws.Cells.Item(46, 1).Value = 'Kitkat 25'
$ws.Cells.Item(46, 3).Value = 25
$ws.Cells.Item(46, 4).Value = 'Kitkat 25.jpg'
$ws.Cells.Item(47, 1).Value = 'Munch 10'
$ws.Cells.Item(47, 3).Value = 10
$ws.Cells.Item(47, 4).Value = 'Munch 10.png'
$ws.Cells.Item(48, 1).Value = 'Sprit 20 rs'
$ws.Cells.Item(48, 3).Value = 10
$ws.Cells.Item(48, 4).Value = 'Sprit 20 rs.jpg'
$ws.Cells.Item(49, 1).Value = 'Fanta20'
$ws.Cells.Item(49, 3).Value = 10
$ws.Cells.Item(49, 4).Value = 'Fanta20.jpg'
$ws.Cells.Item(50, 1).Value = 'Maza10'
$ws.Cells.Item(50, 3).Value = 10
$ws.Cells.Item(50, 4).Value = 'Maza-Tetra-pack.jpg'
$ws.Cells.Item(51, 1).Value = 'Thumsup20'
$ws.Cells.Item(51, 3).Value = 20
$ws.Cells.Item(51, 4).Value = 'Thumsup20.jpg'
$ws.Cells.Item(52, 1).Value = 'Frooti20'
$ws.Cells.Item(52, 3).Value = 20
$ws.Cells.Item(52, 4).Value = 'Frooti20.jpg'
$ws.Cells.Item(53, 1).Value = 'Thumsup25'
$ws.Cells.Item(53, 3).Value = 25
$ws.Cells.Item(53, 4).Value = 'Thusmup Can.jpg'
$ws.Cells.Item(54, 1).Value = 'Dite Coke 25 rs'
$ws.Cells.Item(54, 3).Value = 25
$ws.Cells.Item(54, 4).Value = 'Dite Coke 25 rs.jpg'
$ws.Cells.Item(55, 1).Value = 'Thums up 35'
$ws.Cells.Item(55, 3).Value = 35
$ws.Cells.Item(55, 4).Value = 'Thums up 35.jpg'

# Restore the view: scroll/selection as captured in the saved workbook
[void]$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B46").Select()
